# Stricter separation of departement specific pre mid and post mid courses
#
# This script re-applies a room/section reassignment pass across the
# generated timetable workbook: several rooms used for EC351 / EC301 /
# HS351 / EC306 (Lab) and a handful of elective-basket sessions are
# swapped for different physical rooms, the Room_Allocation /
# Classroom_Utilization / Verification / Basket_Course_Allocations
# summary sheets are refreshed to match, the Executive_Summary
# generation timestamp is bumped, and a department-mismatched
# (CSE-owned HS101) row is dropped from Course_Summary so that
# pre-mid/post-mid course separation is department-scoped.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Timetable + Section_A: identical grids, same edits on both
# ---------------------------------------------------------------------
foreach ($sheetName in @("Timetable", "Section_A")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("E3").Value = "EC351 [C002]"
    $ws.Range("C4").Value = "EC301 [L408]"
    $ws.Range("E4").Value = "HS351 [C202]"
    $ws.Range("C8").Value = "EC306 (Lab) [L106]"
    $ws.Range("B9").Value = "EC351 (Tutorial) [C303]"
    $ws.Range("C9").Value = "EC306 (Lab) [L106]"
}

# ---------------------------------------------------------------------
# Verification: room list summary per course
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Verification")
$ws.Range("I2").Value = "C004, L106"
$ws.Range("I3").Value = "C004, L408"
$ws.Range("I6").Value = "C004, C303, C002"
$ws.Range("I8").Value = "C004, C202"

# ---------------------------------------------------------------------
# Room_Allocation: per-room rollup (room numbers themselves move)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Room_Allocation")

$ws.Range("H2").Value = "MINOR: Cybersecurity, MINOR: Generative Ai, MINOR: Design..."
$ws.Range("H3").Value = "EC351"
$ws.Range("H4").Value = "EC351, EC301 (Tutorial), HS351..."

# row 5 was room C205, now room C202
$ws.Range("A5").Value = "C202"
$ws.Range("D5").Value = "Projector"

# row 6 was room L207, now room C303
$ws.Range("A6").Value = "C303"
$ws.Range("B6").Value = "classroom"
$ws.Range("C6").Value = "96"
$ws.Range("D6").Value = "TV"
$ws.Range("E6").Value = 1
$ws.Range("H6").Value = "EC351 (Tutorial)"
$ws.Range("I6").Value = "0.2"

# row 7 was room L407, now room L106
$ws.Range("A7").Value = "L106"
$ws.Range("B7").Value = "Software Lab"
$ws.Range("C7").Value = "40"
$ws.Range("D7").Value = "Computers"
$ws.Range("E7").Value = 2
$ws.Range("H7").Value = "EC306 (Lab)"
$ws.Range("I7").Value = "0.4"

# row 8 stays room L408, sample course changes
$ws.Range("H8").Value = "EC301"

# ---------------------------------------------------------------------
# Classroom_Allocation: per-session detail rows (same room reshuffle,
# plus matching room-type/capacity/facility lookups + a few basket
# sessions that move to different rooms of a different type)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Classroom_Allocation")

$ws.Range("M6").Value = "C303"

$ws.Range("M10").Value = "L408"

$ws.Range("M11").Value = "L106"

$ws.Range("M12").Value = "L106"

$ws.Range("G19").Value = "large classroom"
$ws.Range("H19").Value = "120"
$ws.Range("I19").Value = "Projector"
$ws.Range("M19").Value = "C002"

$ws.Range("I20").Value = "Projector"
$ws.Range("M20").Value = "C202"

$ws.Range("G25").Value = "classroom"
$ws.Range("H25").Value = "96"
$ws.Range("M25").Value = "C104"

$ws.Range("I26").Value = "TV"
$ws.Range("M26").Value = "C203"

$ws.Range("I27").Value = "TV"
$ws.Range("M27").Value = "C204"

$ws.Range("I28").Value = "TV"
$ws.Range("M28").Value = "C205"

$ws.Range("G30").Value = "classroom"
$ws.Range("H30").Value = "96"
$ws.Range("I30").Value = "Projector"
$ws.Range("M30").Value = "C102"

$ws.Range("G31").Value = "classroom"
$ws.Range("H31").Value = "96"
$ws.Range("I31").Value = "Projector"
$ws.Range("M31").Value = "C104"

$ws.Range("G32").Value = "classroom"
$ws.Range("H32").Value = "96"
$ws.Range("M32").Value = "C202"

$ws.Range("I33").Value = "TV"
$ws.Range("M33").Value = "C203"

$ws.Range("G34").Value = "classroom"
$ws.Range("H34").Value = "96"
$ws.Range("I34").Value = "Projector"
$ws.Range("M34").Value = "C102"

$ws.Range("G35").Value = "classroom"
$ws.Range("H35").Value = "96"
$ws.Range("I35").Value = "Projector"
$ws.Range("M35").Value = "C104"

$ws.Range("G36").Value = "classroom"
$ws.Range("H36").Value = "96"
$ws.Range("M36").Value = "C202"

$ws.Range("I37").Value = "TV"
$ws.Range("M37").Value = "C203"

$ws.Range("G39").Value = "classroom"
$ws.Range("H39").Value = "96"
$ws.Range("M39").Value = "C104"

$ws.Range("I40").Value = "TV"
$ws.Range("M40").Value = "C203"

$ws.Range("I41").Value = "TV"
$ws.Range("M41").Value = "C204"

$ws.Range("I42").Value = "TV"
$ws.Range("M42").Value = "C205"

$ws.Range("G44").Value = "classroom"
$ws.Range("H44").Value = "96"
$ws.Range("I44").Value = "Projector"
$ws.Range("M44").Value = "C102"

$ws.Range("G45").Value = "classroom"
$ws.Range("H45").Value = "96"
$ws.Range("I45").Value = "Projector"
$ws.Range("M45").Value = "C104"

$ws.Range("G46").Value = "classroom"
$ws.Range("H46").Value = "96"
$ws.Range("M46").Value = "C202"

$ws.Range("I47").Value = "TV"
$ws.Range("M47").Value = "C203"

$ws.Range("G49").Value = "classroom"
$ws.Range("H49").Value = "96"
$ws.Range("M49").Value = "C104"

$ws.Range("I50").Value = "TV"
$ws.Range("M50").Value = "C203"

$ws.Range("I51").Value = "TV"
$ws.Range("M51").Value = "C204"

$ws.Range("I52").Value = "TV"
$ws.Range("M52").Value = "C205"

# ---------------------------------------------------------------------
# Executive_Summary: generation timestamp
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Executive_Summary")
$ws.Range("C3").Value = "2026-01-26 12:46"

# ---------------------------------------------------------------------
# Course_Summary: drop the CSE-owned HS101 half-semester row (row 16)
# so ECE pre-mid/post-mid separation no longer leaks other departments'
# half-semester courses into this sheet.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Course_Summary")
$ws.Rows.Item(16).Delete()

# ---------------------------------------------------------------------
# Classroom_Utilization: per-room weekly/daily hour totals follow the
# same room reshuffle (old rooms freed, new rooms now carry the load)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Classroom_Utilization")

# C002
$ws.Range("D3").Value = 1.5
$ws.Range("E3").Value = 0.3
$ws.Range("G3").Value = 3.75

# L106
$ws.Range("D11").Value = 2.5
$ws.Range("E11").Value = 0.5
$ws.Range("G11").Value = 6.25

# C202
$ws.Range("D14").Value = 1.5
$ws.Range("E14").Value = 0.3
$ws.Range("G14").Value = 3.75

# C205
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("G17").Value = 0

# L207
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0
$ws.Range("G19").Value = 0

# C303
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 0.2
$ws.Range("G24").Value = 2.5

# L407
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0
$ws.Range("G36").Value = 0

# ---------------------------------------------------------------------
# Basket_Course_Allocations: elective-basket room assignments move too
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Basket_Course_Allocations")

$ws.Range("C5").Value = "C102"
$ws.Range("C6").Value = "C104"
$ws.Range("C7").Value = "C202"
$ws.Range("C8").Value = "C203"
$ws.Range("C10").Value = "C104"
$ws.Range("C11").Value = "C204"
$ws.Range("C12").Value = "C102"
$ws.Range("C13").Value = "C104"
$ws.Range("C14").Value = "C202"
$ws.Range("C15").Value = "C203"
$ws.Range("C17").Value = "C104"
$ws.Range("C18").Value = "C203"
$ws.Range("C19").Value = "C204"
$ws.Range("C20").Value = "C205"
